$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3160.1428
$ws.Range("I74").Value = 3056.1875
$ws.Range("K74").Value = 3056.1875
$ws.Range("M74").Value = -2120.1875
$ws.Range("H77").Value = 3160.1428
$ws.Range("I77").Value = 3056.1875
$ws.Range("K77").Value = 15280.9375
$ws.Range("M77").Value = -10600.9375
$ws.Range("H86").Value = 3008.16
$ws.Range("I86").Value = 6000
$ws.Range("J86").Value = 2748
$ws.Range("K86").Value = 6000
$ws.Range("L86").Value = 2748
$ws.Range("M86").Value = -4877
$ws.Range("N86").Value = -4994
$ws.Range("H89").Value = 3008.16
$ws.Range("I89").Value = 6000
$ws.Range("J89").Value = 2748
$ws.Range("K89").Value = 30000
$ws.Range("L89").Value = 13740
$ws.Range("M89").Value = -24384
$ws.Range("N89").Value = -24972

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3071.0625
$ws.Range("I61").Value = 2163.8462
$ws.Range("J61").Value = 7002.3335
$ws.Range("K61").Value = 2163.8462
$ws.Range("L61").Value = 7002.3335
$ws.Range("M61").Value = -1951.8462
$ws.Range("N61").Value = -7426.3335
$ws.Range("H97").Value = 1686.6945
$ws.Range("I97").Value = 1297.0333
$ws.Range("J97").Value = 3635
$ws.Range("K97").Value = 1297.0333
$ws.Range("L97").Value = 3635
$ws.Range("M97").Value = -801.0333000000001
$ws.Range("N97").Value = -4627
$ws.Range("H136").Value = 3071.0625
$ws.Range("I136").Value = 2163.8462
$ws.Range("J136").Value = 7002.3335
$ws.Range("K136").Value = 6491.5386
$ws.Range("L136").Value = 21007.0005
$ws.Range("M136").Value = -3941.5386
$ws.Range("N136").Value = -26107.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2467.6667
$ws.Range("I20").Value = 2412
$ws.Range("K20").Value = 2412
$ws.Range("M20").Value = -2165
$ws.Range("H86").Value = 4336.8335
$ws.Range("I86").Value = 3633.8333
$ws.Range("J86").Value = 5742.8335
$ws.Range("K86").Value = 3633.8333
$ws.Range("L86").Value = 5742.8335
$ws.Range("M86").Value = -2510.8333
$ws.Range("N86").Value = -7988.8335
$ws.Range("H89").Value = 4336.8335
$ws.Range("I89").Value = 3633.8333
$ws.Range("J89").Value = 5742.8335
$ws.Range("K89").Value = 18169.1665
$ws.Range("L89").Value = 28714.1675
$ws.Range("M89").Value = -12553.1665
$ws.Range("N89").Value = -39946.1675

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2869.125
$ws.Range("I31").Value = 2441.276
$ws.Range("J31").Value = 3997.0908
$ws.Range("K31").Value = 2441.276
$ws.Range("L31").Value = 3997.0908
$ws.Range("M31").Value = -2146.276
$ws.Range("N31").Value = -4587.0908
$ws.Range("H34").Value = 2869.125
$ws.Range("I34").Value = 2441.276
$ws.Range("J34").Value = 3997.0908
$ws.Range("K34").Value = 2441.276
$ws.Range("L34").Value = 3997.0908
$ws.Range("M34").Value = -2239.276
$ws.Range("N34").Value = -4401.0908
$ws.Range("H58").Value = 1951.3871
$ws.Range("I58").Value = 1888.6296
$ws.Range("J58").Value = 2375
$ws.Range("K58").Value = 1888.6296
$ws.Range("L58").Value = 2375
$ws.Range("M58").Value = -1685.6296
$ws.Range("N58").Value = -2781
$ws.Range("H132").Value = 1223.8422
$ws.Range("I132").Value = 1090.129
$ws.Range("J132").Value = 1816
$ws.Range("K132").Value = 3270.387
$ws.Range("L132").Value = 5448
$ws.Range("M132").Value = -740.3869999999997
$ws.Range("N132").Value = -10508
$ws.Range("H134").Value = 1369.8064
$ws.Range("I134").Value = 1364.9524
$ws.Range("J134").Value = 1380
$ws.Range("K134").Value = 4094.857199999999
$ws.Range("L134").Value = 4140
$ws.Range("M134").Value = -1559.857199999999
$ws.Range("N134").Value = -9210
$ws.Range("H136").Value = 1951.3871
$ws.Range("I136").Value = 1888.6296
$ws.Range("J136").Value = 2375
$ws.Range("K136").Value = 5665.8888
$ws.Range("L136").Value = 7125
$ws.Range("M136").Value = -3115.8888
$ws.Range("N136").Value = -12225

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H59").Value = 5600

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 12611.5
$ws.Range("I63").Value = 10103
$ws.Range("J63").Value = 15120
$ws.Range("K63").Value = 10103
$ws.Range("L63").Value = 15120
$ws.Range("N63").Value = -16492
$ws.Range("M63").Value = -9417
$ws.Range("H66").Value = 12611.5
$ws.Range("I66").Value = 10103
$ws.Range("J66").Value = 15120
$ws.Range("K66").Value = 30309
$ws.Range("L66").Value = 45360
$ws.Range("N66").Value = -52224
$ws.Range("M66").Value = -26877
$ws.Range("H80").Value = 2277.1304
$ws.Range("I80").Value = 2017.5
$ws.Range("J80").Value = 2476.8462
$ws.Range("K80").Value = 2017.5
$ws.Range("L80").Value = 2476.8462
$ws.Range("M80").Value = -1019.5
$ws.Range("N80").Value = -4472.8462
$ws.Range("H83").Value = 2277.1304
$ws.Range("I83").Value = 2017.5
$ws.Range("J83").Value = 2476.8462
$ws.Range("K83").Value = 10087.5
$ws.Range("L83").Value = 12384.231
$ws.Range("M83").Value = -5095.5
$ws.Range("N83").Value = -22368.231

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3253.28
$ws.Range("I40").Value = 2906.35
$ws.Range("J40").Value = 4641
$ws.Range("K40").Value = 2906.35
$ws.Range("L40").Value = 4641
$ws.Range("M40").Value = -2770.35
$ws.Range("N40").Value = -4913
$ws.Range("H132").Value = 3415.4546
$ws.Range("I132").Value = 1480.3273
$ws.Range("J132").Value = 13091.091
$ws.Range("K132").Value = 4440.9819
$ws.Range("L132").Value = 39273.273
$ws.Range("M132").Value = -1910.9819
$ws.Range("N132").Value = -44333.273

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 17226
$ws.Range("I14").Value = 2968
$ws.Range("J14").Value = 60000
$ws.Range("K14").Value = 2968
$ws.Range("L14").Value = 60000
$ws.Range("M14").Value = -2800
$ws.Range("N14").Value = -60336
$ws.Range("H81").Value = 640.5925999999999
$ws.Range("I81").Value = 623.73914
$ws.Range("J81").Value = 737.5
$ws.Range("K81").Value = 1247.47828
$ws.Range("L81").Value = 1475
$ws.Range("M81").Value = -186.47828
$ws.Range("N81").Value = -3597
$ws.Range("H84").Value = 640.5925999999999
$ws.Range("I84").Value = 623.73914
$ws.Range("J84").Value = 737.5
$ws.Range("K84").Value = 6237.3914
$ws.Range("L84").Value = 7375
$ws.Range("M84").Value = -933.3914000000004
$ws.Range("N84").Value = -17983
$ws.Range("H100").Value = 536.86664
$ws.Range("I100").Value = 472.72726
$ws.Range("J100").Value = 713.25
$ws.Range("K100").Value = 945.45452
$ws.Range("L100").Value = 1426.5
$ws.Range("M100").Value = -404.45452
$ws.Range("N100").Value = -2508.5
$ws.Range("H136").Value = 741.775
$ws.Range("I136").Value = 678.9375
$ws.Range("J136").Value = 993.125
$ws.Range("K136").Value = 2036.8125
$ws.Range("L136").Value = 2979.375
$ws.Range("M136").Value = 513.1875
$ws.Range("N136").Value = -8079.375
